$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.1958502231010186
$ws.Range("D2").Value = 0.1172342349564133
$ws.Range("E2").Value = 0.06542532866897588
$ws.Range("F2").Value = 8.797736409582797
$ws.Range("G2").Value = 0.002621317206607812
$ws.Range("I2").Value = 5.987256124704885
$ws.Range("J2").Value = 0.2975686261519286
$ws.Range("N2").Value = 2.873581687696401
$ws.Range("C3").Value = 0.1894650465831802
$ws.Range("D3").Value = 0.1015593299189135
$ws.Range("E3").Value = 0.06589401896876534
$ws.Range("F3").Value = 8.431914133496718
$ws.Range("G3").Value = 0.002639233995562918
$ws.Range("I3").Value = 5.707567514089959
$ws.Range("J3").Value = 0.289671771849811
$ws.Range("N3").Value = 2.562605684679681
$ws.Range("C4").Value = 0.1857308733923873
$ws.Range("D4").Value = 0.09201031273020988
$ws.Range("E4").Value = 0.06621019150032392
$ws.Range("F4").Value = 8.215872885662009
$ws.Range("G4").Value = 0.002650734123936161
$ws.Range("I4").Value = 5.541798128587914
$ws.Range("J4").Value = 0.2850930252612613
$ws.Range("N4").Value = 2.371325805375818
$ws.Range("C5").Value = 0.1842546618158849
$ws.Range("D5").Value = 0.08813564783228856
$ws.Range("E5").Value = 0.06634616491355505
$ws.Range("F5").Value = 8.129920768793397
$ws.Range("G5").Value = 0.002655547100376327
$ws.Range("I5").Value = 5.475692672114945
$ws.Range("J5").Value = 0.2832932862343966
$ws.Range("N5").Value = 2.293303068605894
$ws.Range("C6").Value = 0.1840122468715606
$ws.Range("D6").Value = 0.08749320094040058
$ws.Range("E6").Value = 0.06636917360730088
$ws.Range("F6").Value = 8.115772511590137
$ws.Range("G6").Value = 0.002656353967446463
$ws.Range("I6").Value = 5.464801821958702
$ws.Range("J6").Value = 0.2829983857940022
$ws.Range("N6").Value = 2.2803432614038
$ws.Range("C7").Value = 0.1857107822456214
$ws.Range("D7").Value = 0.09195799330225896
$ws.Range("E7").Value = 0.06621199642550035
$ws.Range("F7").Value = 8.214705347696054
$ws.Range("G7").Value = 0.002650798519678107
$ws.Range("I7").Value = 5.540900815911414
$ws.Range("J7").Value = 0.2850684878173695
$ws.Range("N7").Value = 2.370273851392596
$ws.Range("C8").Value = 0.1936090963418451
$ws.Range("D8").Value = 0.1118124364932385
$ws.Range("E8").Value = 0.06558103358789324
$ws.Range("F8").Value = 8.669778124637389
$ws.Range("G8").Value = 0.002627391967533556
$ws.Range("I8").Value = 5.88955041679651
$ws.Range("J8").Value = 0.2947886871178582
$ws.Range("N8").Value = 2.766433886209882
$ws.Range("C9").Value = 0.2106414648431212
$ws.Range("D9").Value = 0.1514564176480349
$ws.Range("E9").Value = 0.06456952801849525
$ws.Range("F9").Value = 9.63354846943281
$ws.Range("G9").Value = 0.002585402634400355
$ws.Range("I9").Value = 6.623023479938013
$ws.Range("J9").Value = 0.3160753154067066
$ws.Range("N9").Value = 3.540180268007646
$ws.Range("C10").Value = 0.2241923068920357
$ws.Range("D10").Value = 0.1811780494233233
$ws.Range("E10").Value = 0.06396486864055362
$ws.Range("F10").Value = 10.39005241784611
$ws.Range("G10").Value = 0.0025568659430434
$ws.Range("I10").Value = 7.195896567551074
$ws.Range("J10").Value = 0.333194154190366
$ws.Range("N10").Value = 4.10623028343673
$ws.Range("C11").Value = 0.2306045628622257
$ws.Range("D11").Value = 0.1948661997597583
$ws.Range("E11").Value = 0.0637200846189323
$ws.Range("F11").Value = 10.74586501054262
$ws.Range("G11").Value = 0.002544370028750118
$ws.Range("I11").Value = 7.464742341545332
$ws.Range("J11").Value = 0.3413319376062276
$ws.Range("N11").Value = 4.363110593465422
$ws.Range("C12").Value = 0.2330703087452832
$ws.Range("D12").Value = 0.2000769047677977
$ws.Range("E12").Value = 0.06363176626004829
$ws.Range("F12").Value = 10.88237983455616
$ws.Range("G12").Value = 0.002539706637760836
$ws.Range("I12").Value = 7.567806037219441
$ws.Range("J12").Value = 0.344466369326085
$ws.Range("N12").Value = 4.460285735713398
$ws.Range("C13").Value = 0.2325375639518938
$ws.Range("D13").Value = 0.1989534193569966
$ws.Range("E13").Value = 0.06365059219857017
$ws.Range("F13").Value = 10.85289831267175
$ws.Range("G13").Value = 0.00254070795532579
$ws.Range("I13").Value = 7.545552298190216
$ws.Range("J13").Value = 0.3437889232332338
$ws.Range("N13").Value = 4.439361943450422
$ws.Range("C14").Value = 0.2308066570763003
$ws.Range("D14").Value = 0.1952943215653136
$ws.Range("E14").Value = 0.06371273077781581
$ws.Range("F14").Value = 10.75705998972887
$ws.Range("G14").Value = 0.00254398500270566
$ws.Range("I14").Value = 7.473195820422347
$ws.Range("J14").Value = 0.3415887358662673
$ws.Range("N14").Value = 4.371107314139238
$ws.Range("C15").Value = 0.2297513772034563
$ws.Range("D15").Value = 0.1930566714435145
$ws.Range("E15").Value = 0.06375136303320694
$ws.Range("F15").Value = 10.69859056311242
$ws.Range("G15").Value = 0.00254600117601878
$ws.Range("I15").Value = 7.429041384337552
$ws.Range("J15").Value = 0.3402480126896279
$ws.Range("N15").Value = 4.329286057409945
$ws.Range("C16").Value = 0.2237784006953518
$ws.Range("D16").Value = 0.1802871389869551
$ws.Range("E16").Value = 0.06398147696506307
$ws.Range("F16").Value = 10.36704236982945
$ws.Range("G16").Value = 0.002557692242617177
$ws.Range("I16").Value = 7.178498824053293
$ws.Range("J16").Value = 0.3326695870509297
$ws.Range("N16").Value = 4.089429168003846
$ws.Range("C17").Value = 0.2201789975729866
$ws.Range("D17").Value = 0.1724985986842285
$ws.Range("E17").Value = 0.06413041443747147
$ws.Range("F17").Value = 10.1667062343887
$ws.Range("G17").Value = 0.002564987800613109
$ws.Range("I17").Value = 7.026960502195379
$ws.Range("J17").Value = 0.3281118998073822
$ws.Range("N17").Value = 3.94211849063862
$ws.Range("C18").Value = 0.218131929739144
$ws.Range("D18").Value = 0.1680344134783525
$ws.Range("E18").Value = 0.06421892784428351
$ws.Range("F18").Value = 10.05257096379484
$ws.Range("G18").Value = 0.002569229796185713
$ws.Range("I18").Value = 6.940570972344062
$ws.Range("J18").Value = 0.3255232593855055
$ws.Range("N18").Value = 3.857331695637072
$ws.Range("C19").Value = 0.217442758104113
$ws.Range("D19").Value = 0.1665254978751705
$ws.Range("E19").Value = 0.06424938556293291
$ws.Range("F19").Value = 10.01411152534195
$ws.Range("G19").Value = 0.002570673964643469
$ws.Range("I19").Value = 6.911451295156326
$ws.Range("J19").Value = 0.3246523547284568
$ws.Range("N19").Value = 3.828614786363971
$ws.Range("C20").Value = 0.2205597441278826
$ws.Range("D20").Value = 0.1733260673151165
$ws.Range("E20").Value = 0.06411426482689109
$ws.Range("F20").Value = 10.18791857106834
$ws.Range("G20").Value = 0.002564206446589612
$ws.Range("I20").Value = 7.043011676941944
$ws.Range("J20").Value = 0.3285936592506431
$ws.Range("N20").Value = 3.95780600327754
$ws.Range("C21").Value = 0.2313140311191262
$ws.Range("D21").Value = 0.1963683200737592
$ws.Range("E21").Value = 0.0636943602197011
$ws.Range("F21").Value = 10.78516101602702
$ws.Range("G21").Value = 0.002543020604626201
$ws.Range("I21").Value = 7.494413943779364
$ws.Range("J21").Value = 0.3422335302183086
$ws.Range("N21").Value = 4.391158149571083
$ws.Range("C22").Value = 0.2385624967866136
$ws.Range("D22").Value = 0.211588494353407
$ws.Range("E22").Value = 0.06344544705117094
$ws.Range("F22").Value = 11.18589639511981
$ws.Range("G22").Value = 0.002529573310048424
$ws.Range("I22").Value = 7.796798367749318
$ws.Range("J22").Value = 0.3514571454263091
$ws.Range("N22").Value = 4.67379181795809
$ws.Range("C23").Value = 0.2346730706499045
$ws.Range("D23").Value = 0.2034493933969372
$ws.Range("E23").Value = 0.06357595385189718
$ws.Range("F23").Value = 10.9710308850299
$ws.Range("G23").Value = 0.002536714323400574
$ws.Range("I23").Value = 7.634711082091201
$ws.Range("J23").Value = 0.3465051926165188
$ws.Range("N23").Value = 4.523002190001307
$ws.Range("C24").Value = 0.2203875393657597
$ws.Range("D24").Value = 0.1729519267805983
$ws.Range("E24").Value = 0.06412155707434586
$ws.Range("F24").Value = 10.17832523959248
$ws.Range("G24").Value = 0.002564559548160407
$ws.Range("I24").Value = 7.035752666648989
$ws.Range("J24").Value = 0.3283757574548503
$ws.Range("N24").Value = 3.950713976768498
$ws.Range("C25").Value = 0.2058588626372284
$ws.Range("D25").Value = 0.1406393807135657
$ws.Range("E25").Value = 0.06481891287796948
$ws.Range("F25").Value = 9.364689707039133
$ws.Range("G25").Value = 0.0025963504182308
$ws.Range("I25").Value = 6.418905353816967
$ws.Range("J25").Value = 0.3100661645182754
$ws.Range("N25").Value = 3.331249627311138
